$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: a known-default-styled cell to copy style from when we need to force text storage
# for numeric-looking strings (Excel would otherwise auto-convert them to numbers).
$plainStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = '54.221.27'
$ws.Range("E2").Value = '  +1.17%  '

$ws.Range("D3").Value = '2.289.89'
$ws.Range("E3").Value = '  +3.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '495.81'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +2.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.03'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  +1.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  +2.46%  '

$ws.Range("D9").Value = '2.289.05'
$ws.Range("E9").Value = '  +2.94%  '

$ws.Range("E10").Value = '  +4.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  +2.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.326'
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = '  +4.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.64'
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = '  -0.71%  '

$ws.Range("D14").Value = '2.696.38'
$ws.Range("E14").Value = '  +3.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.80'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +3.71%  '

$ws.Range("D16").Value = '54.268.85'
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("E17").Value = '  +1.67%  '

$ws.Range("D18").Value = '2.303.68'
$ws.Range("E18").Value = '  +3.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.04'
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  +5.37%  '

$ws.Range("E20").Value = '  +4.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '301.36'
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = '  +1.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.44'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  +5.83%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.37'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  -2.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.48'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  -1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("E27").Value = '  +2.94%  '

$ws.Range("D28").Value = '2.394.10'
$ws.Range("E28").Value = '  +3.07%  '

$ws.Range("E29").Value = '  +4.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.05'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.70'
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  +0.53%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.60'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  +1.98%  '

$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '0.0₃0689'
$ws.Range("E33").Value = '  +2.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  +3.41%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("E37").Value = '  +2.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.72'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  +2.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.909'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  +9.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.19'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +4.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.70'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  +4.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.58'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.40'
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = '  +3.26%  '

$ws.Range("E44").Value = '  +2.90%  '

$ws.Range("E45").Value = '  +3.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '127.37'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  +4.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.77'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  +3.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0888'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  +1.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.548'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  +3.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '238.64'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("E51").Value = '  +3.67%  '
